# Update "想去人数" (interest count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 8835
$ws1.Range("F10").Value = 816
$ws1.Range("F11").Value = 333
$ws1.Range("F18").Value = 283
$ws1.Range("F21").Value = 1073

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 8835
$ws4.Range("F12").Value = 816
$ws4.Range("F13").Value = 333
$ws4.Range("F20").Value = 283
$ws4.Range("F23").Value = 1073
